$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the census-note row (old row 2) entirely, shifting remaining rows up.
$ws.Rows.Item(2).Delete()

# Remove the 1989/2002 year columns, keeping only the 2014 figures
# (old column B and C), shifting column D left into column B.
$ws.Range("B:C").EntireColumn.Delete()

# Select A2 to match the saved selection state.
$ws.Range("A2").Select()

# Rename the worksheet/tab from "1" to the municipality name.
$ws.Name = "ნინოწმინდა"
